# Apply the "add the missing repo link in pptx" edit to slide 2's
# "Text Box 147" shape:
#   1) Split the empty "()" run into "(" + hyperlinked URL + ")".
#   2) Resize/reposition the shape (it widens to fit the new text).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Locate the "Text Box 147" shape (holds the "...()..." sentence).
$shp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $cand = $s.Shapes.Item($i)
    if ($cand.Name -eq "Text Box 147") {
        $shp = $cand
    }
}

# --- 1) Split the "()" run and insert the hyperlinked URL ----------------
$tf = $shp.TextFrame
$tr = $tf.TextRange
$fullText = $tr.Text
$parenIdx = $fullText.IndexOf("()")

$url = "https://github.com/wszqkzqk/pypvz"

# Replace "()" with "(<url>)" in-place (keeps the run's existing formatting).
$parenRange = $tr.Characters($parenIdx + 1, 2)
$parenRange.Text = "(" + $url + ")"

# Re-locate the URL text and turn it into its own hyperlinked run.
$fullText2 = $tr.Text
$urlIdx = $fullText2.IndexOf($url)
$urlRange = $tr.Characters($urlIdx + 1, $url.Length)
$urlRange.ActionSettings.Item(1).Hyperlink.Address = $url

# --- 2) Reposition / resize the shape (done last: the textbox has
#        spAutoFit, so its Height is recomputed on every text/width
#        change -- setting the geometry after editing the text makes
#        the explicit values below stick). ------------------------------
# Target EMUs: off x=377018 y=1288549, ext cx=11591170 cy=5013039
# (literal point values below are pre-computed so that the host's
# float32 + floor EMU conversion lands exactly on the target EMU values)
$shp.Left = 29.686457692913386
$shp.Top = 101.46055228110237
$shp.Width = 912.6905511811024
$shp.Height = 394.72749331496067
